$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44181
$ws.Range("M2").Value = 220
$ws.Range("N2").Value = 17000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 17000
$ws.Range("Q2").Value = '$/caja 18 kilos'
$ws.Range("R2").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S2").Value = 944
$ws.Range("T2").Value = 18
$ws.Range("D3").Value = 44174
$ws.Range("L3").Value = 'Especial'
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = '$/caja 10 kilos'
$ws.Range("R3").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S3").Value = 1500
$ws.Range("T3").Value = 10
$ws.Range("D4").Value = 44186
$ws.Range("K4").Value = 'Dina'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 150
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 833
$ws.Range("T4").Value = 18
$ws.Range("D6").Value = 44179
$ws.Range("M6").Value = 150
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("S6").Value = 1000
$ws.Range("D7").Value = 44161
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 150
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Provincia de Limarí'
$ws.Range("S7").Value = 1111
$ws.Range("T7").Value = 18
$ws.Range("D8").Value = 44167
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = '$/caja 16 kilos granel'
$ws.Range("R8").Value = 'Provincia de Limarí'
$ws.Range("S8").Value = 938
$ws.Range("T8").Value = 16
$ws.Range("D9").Value = 44172
$ws.Range("M9").Value = 120
$ws.Range("D10").Value = 44172
$ws.Range("K10").Value = 'Castle Brite'
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 11000
$ws.Range("O10").Value = 11000
$ws.Range("P10").Value = 11000
$ws.Range("Q10").Value = '$/caja 10 kilos'
$ws.Range("R10").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S10").Value = 1100
$ws.Range("T10").Value = 10
$ws.Range("D11").Value = 44187
$ws.Range("M11").Value = 120
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 16000
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 889
$ws.Range("D12").Value = 44168
$ws.Range("K12").Value = 'Castle Brite'
$ws.Range("M12").Value = 250
$ws.Range("N12").Value = 10000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 10000
$ws.Range("Q12").Value = '$/caja 10 kilos'
$ws.Range("T12").Value = 10
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 17000
$ws.Range("O13").Value = 17000
$ws.Range("P13").Value = 17000
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("R13").Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 944
$ws.Range("T13").Value = 18
$ws.Range("D14").Value = 44160
$ws.Range("K14").Value = 'Dina'
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 20000
$ws.Range("Q14").Value = '$/caja 15 kilos'
$ws.Range("S14").Value = 1333
$ws.Range("T14").Value = 15
$ws.Range("D15").Value = 44162
$ws.Range("M15").Value = 200
$ws.Range("Q15").Value = '$/caja 16 kilos granel'
$ws.Range("S15").Value = 1062
$ws.Range("T15").Value = 16
$ws.Range("D16").Value = 44162
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("Q16").Value = '$/caja 16 kilos granel'
$ws.Range("S16").Value = 938
$ws.Range("T16").Value = 16
$ws.Range("D17").Value = 44176
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 17000
$ws.Range("P17").Value = 17000
$ws.Range("Q17").Value = '$/caja 18 kilos granel'
$ws.Range("S17").Value = 944
$ws.Range("T17").Value = 18
